$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Queries")

$sql1 = @'
SELECT M.AgentID as [Agent ID],A.AgentName as[Agent Name],A.TeamName as [Team Name],A.SupervisorName as[Supervisor Name],sum(ACDCalls) AS [Total Interaction],[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])) as [Total Interaction Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])/nullif(SUM(ACDCalls),0))AS [Avg Interaction Time],SUM([TotalChat]) as [Total Chat],
[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])) as[Total Chat Time],[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])/nullif(sum([TotalChat]),0))AS [Avg Chat Time],SUM([TotalAudioIP]) as [Total Audio IP],[dbo].[SECONDSTOhhmmss](sum([TotalAudioIPTime])) [Total AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalAudioIPTime])/nullif(SUM([TotalAudioIP]),0))AS [Avg AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalInteractionTime)/nullif(SUM(ACDCalls),0))AS [Avg Talk Time],[dbo].[SECONDSTOhhmmss](SUM(TotalAfterCallTime)) as [Total After Call Time],[dbo].[SECONDSTOhhmmss](SUM(TotalAvailTime)) as [Total Avail Time],[dbo].[SECONDSTOhhmmss](SUM(TotalAuxTime)) as [Total AUX Time],SUM(ExtensionCalls) AS [Extension Interaction],[dbo].[SECONDSTOhhmmss](sum([TotalExtensionTime])) as [Total Extension Time],
[dbo].[SECONDSTOhhmmss](sum(TotalExtensionTime)/nullif(sum(ExtensionCalls),0))AS [Avg Extension Time],[dbo].[SECONDSTOhhmmss](SUM(TotalStaffedTime)) as [Total Time Staffed],[dbo].[SECONDSTOhhmmss](SUM(TotalHoldTime)) as [Total Hold Time] 
FROM [OCM_AgentHistoricalReport] M WITH(NOLOCK)
INNER JOIN fn_AgentHierarchy('na','1','1') A ON  A.[AgentId]=M.[AgentID]
WHERE [ReportDateTime]>='ReportBeforeDate' and [ReportDateTime]<='ReportAfterDate'
GROUP BY M.[AgentID], A.[AgentName],A.[TeamName],A.[SupervisorName]
Order by [Agent Name];
'@

$sql2 = @'
SELECT Dateint AS [Date],AgentID as [Agent ID],ISNULL(A.FirstName,'')+' '+ ISNULL(A.LastName,'') AS [Agent Name],ISNULL(C.TeamName,' ') AS TeamName,
ISNULL(B.FirstName,'NA')+' '+ ISNULL(B.LastName,'') AS SupervisorName,SUM(ACDCalls) AS TotalInteraction,
[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])) [Total Interaction Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalInteractionTime])/nullif(SUM(ACDCalls),0))AS [Avg Interaction Time],
SUM([TotalChat]) as [Total Chat],[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])) as [Total Chat Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalChatTime])/nullif(sum([TotalChat]),0))AS [Avg Chat Time],
SUM([TotalAudioIP]) [Total Audio IP],[dbo].[SECONDSTOhhmmss](sum([TotalAudioIPTime])) [Total AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM([TotalAudioIPTime])/nullif(SUM([TotalAudioIP]),0))AS [Avg AudioIP Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalInteractionTime)/nullif(SUM(ACDCalls),0))AS [Avg Talk Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalAfterCallTime)) as [Total After Call Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalAvailTime)) as [Total Avail Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalAuxTime)) as [Total Aux Time],
SUM(ExtensionCalls) AS [Extension Interaction],
[dbo].[SECONDSTOhhmmss](sum([TotalExtensionTime])) [Total Extension Time],
[dbo].[SECONDSTOhhmmss](sum(TotalExtensionTime)/nullif(sum(ExtensionCalls),0))AS [Avg Extension Time],
[dbo].[SECONDSTOhhmmss](SUM(TotalStaffedTime)) as [Total Time Staffed],
[dbo].[SECONDSTOhhmmss](SUM(TotalHoldTime)) as [Total Hold Time]  
FROM [OCM_AgentHistoricalReport] M WITH(NOLOCK)
LEFT JOIN [AGT_Agent] A WITH(NOLOCK)  ON A.AvayaLoginID = M.[AgentID] 
LEFT JOIN [AGT_Agent] B WITH(NOLOCK)  ON A.[PrimarySupervisorID]=B.ID LEFT JOIN [AGT_Teams] C WITH(NOLOCK) ON C.TeamID = A.TeamID
LEFT JOIN[dbo].[AGT_Teams] P WITH(NOLOCK) ON C.ParentID = P.TeamID 
WHERE [ReportDateTime]>='ReportBeforeDate' AND [ReportDateTime]<='ReportAfterDate' AND [AgentID] LIKE 'AgentIdCapturedFromUI'  
GROUP BY  [Dateint], [AgentID],B.[FirstName],B.[LastName],C.TeamName,A.FirstName,A.LastName ORDER BY [Dateint] ASC;
'@

# F2: first query, wrap-text style only
$ws.Range("F2").Value = $sql1
$ws.Range("F2").WrapText = $true

# G2: second query, wrap-text + centered style
$ws.Range("G2").Value = $sql2
$ws.Range("G2").WrapText = $true
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").VerticalAlignment = -4108

# Row 2 grows to fit the large wrapped text
$ws.Rows(2).RowHeight = 409.5

# Column widths: F best-fits to the long single-line query, G is hand widened
$ws.Range("F1").EntireColumn.ColumnWidth = 82
$ws.Range("G1").EntireColumn.ColumnWidth = 56.67

# Move the selection/view over to the newly added columns
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("G2").Select()
